$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sainsburys sheet: drop stale hyperlinks, record the Fairy liquid qty,
# and tidy up the leftover fill style on the category rows.
# ---------------------------------------------------------------------------
$sains = $wb.Worksheets.Item("Sainsburys")

$sains.Hyperlinks.Delete()

$sains.Range("C8").Value = 2

# Rows 9:22 previously carried a redundant "applyFill" style; restore the
# plain/default look by painting the format from an already-plain cell.
$sains.Range("A2").Copy()
$sains.Range("A9:A22").PasteSpecial(-4122)

$sains.Columns.Item(2).ColumnWidth = 150.99869791666666

$sains.Range("B11").Select() | Out-Null

# ---------------------------------------------------------------------------
# Tesco sheet: populate with the new store-comparison data.
# ---------------------------------------------------------------------------
$tesco = $wb.Worksheets.Item("Tesco")

# Header row: copy the header format from Sainsburys, then fill in the text.
$sains.Range("A1:C1").Copy()
$tesco.Range("A1:C1").PasteSpecial(-4122)
$tesco.Range("A1").Value = "Product Name"
$tesco.Range("B1").Value = "URL"
$tesco.Range("C1").Value = "Price"

# Data rows: stamp the plain body style (copied from a plain Sainsburys
# data cell) across A2:C8 first so blank cells still keep a style ref,
# then fill in the values.
$sains.Range("A2").Copy()
$tesco.Range("A2:C8").PasteSpecial(-4122)

$tesco.Range("A2").Value = "Kingsmill Medium Sliced 50/50 Bread 800g"
$tesco.Range("B2").Value = "https://www.tesco.com/groceries/en-GB/products/261738730"
$tesco.Range("C2").Value = 2.2

$tesco.Range("A3").Value = "Hovis Medium Sliced Wholemeal Bread 800g"
$tesco.Range("B3").Value = "https://www.tesco.com/groceries/en-GB/products/255000362"

$tesco.Range("A4").Value = "Tesco British Semi Skimmed Milk 2.272L 4 Pints"
$tesco.Range("B4").Value = "https://www.tesco.com/groceries/en-GB/products/254656543"

$tesco.Range("A5").Value = "British Whole Milk 2.27L (4 pint)"

$tesco.Range("A6").Value = "Tesco Semi Skimmed Milk 3.408L/6 Pints"
$tesco.Range("B6").Value = "https://www.tesco.com/groceries/en-GB/products/255986260"

$tesco.Range("A7").Value = "British Whole Milk 3.4L (6 pint)"

$tesco.Range("A8").Value = "Fairy Platinum Quickwash Washing Up Liquid, Original 625ml"

$tesco.Columns.Item(1).ColumnWidth = 51.830729166666664
$tesco.Columns.Item(2).ColumnWidth = 53.330729166666664

$tesco.Activate() | Out-Null
$tesco.Range("A14").Select() | Out-Null
